$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff represents a reordering: the data that was in row 13 moves to
# row 14, and the data that was in row 14 moves to row 13 (all columns).
# Swap the two rows cell-by-cell, being careful to preserve the original
# cell types (several "numeric looking" values such as "1", "5", times
# and dates are actually stored as text, not as numbers/dates).

$numericCols = @("A","B","E","Q","R","S")
$boolCols    = @("AD","AE","AG")
$textCols    = @("D","F","G","H","I","J","K","M","P","T","U","V","W","Y","Z","AA","AB","AT","AW","AX","AY")

function Get-RowData($row) {
  $data = @{}
  foreach ($col in $numericCols) { $data[$col] = $ws.Range("$col$row").Value2 }
  foreach ($col in $boolCols)    { $data[$col] = $ws.Range("$col$row").Value2 }
  foreach ($col in $textCols)    { $data[$col] = $ws.Range("$col$row").Text }
  return $data
}

$row13 = Get-RowData 13
$row14 = Get-RowData 14

function Set-RowData($row, $data) {
  foreach ($col in $numericCols) {
    $cell = $ws.Range("$col$row")
    if ($null -ne $data[$col]) { $cell.Value2 = $data[$col] } else { $cell.Value2 = $null }
  }
  foreach ($col in $boolCols) {
    $cell = $ws.Range("$col$row")
    if ($null -ne $data[$col]) { $cell.Value2 = $data[$col] } else { $cell.Value2 = $null }
  }
  foreach ($col in $textCols) {
    $cell = $ws.Range("$col$row")
    $cell.NumberFormat = "@"
    if ($data[$col] -ne "") { $cell.Value2 = $data[$col] } else { $cell.Value2 = $null }
  }
}

Set-RowData 13 $row14
Set-RowData 14 $row13
